$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the inventory table with five more rows (ids 10-14), each
# stamped with the caja-mensual update date in column H (fechaActualizacion).
$startRow = 11
$startId = 10
$fecha = "2023-09-11"

for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $id = $startId + $i

    $ws.Cells.Item($row, 1).Value = $id

    # Force the date-looking string to stay plain text (matches the rest of
    # the sheet's fechaActualizacion column), then drop back to the default
    # "Normal" style so no stray number format sticks to the cell.
    $ws.Cells.Item($row, 8).Value = "'" + $fecha
    $ws.Cells.Item($row, 8).Style = "Normal"
}
